$d = $word.ActiveDocument

# The document body's only paragraph holds a single inline picture
# (the chapter question's image). Remove the picture, leaving the
# paragraph mark in place (an empty paragraph), matching the target
# edit that strips the <w:drawing> run out of the paragraph.
if ($d.InlineShapes.Count -gt 0) {
    for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
        $d.InlineShapes.Item($i).Delete()
    }
}
